$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new translation-table rows right after the existing
# "AsmntAggFun" row (row 32), pushing every row below down by two.
$ws.Rows.Item(33).EntireRow.Insert()
$ws.Rows.Item(34).EntireRow.Insert()

$ws.Range("A33").Value = "AsmntAggPeriod"
$ws.Range("B33").Value = "DA"

$ws.Range("A34").Value = "AsmntAggPeriodUnit"
$ws.Range("B34").Value = "DA"

# Restore the view/selection state recorded for the edited sheet.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select()
